$d = $word.ActiveDocument

# The "Requisitos" bullet list currently has three weak-requirement lines
# (LOM3036, LOM3082, LOM3057). Replace the whole block with a single line
# for LOM3013, keeping the same paragraph/line-break structure.
$found = $d.Content.Find.Execute(
    "LOM3036*LOM3057*Requisito fraco)",
    $false, $false, $true, $false, $false,
    $true, 1, $false,
    "LOM3013 -  Ciência dos Materiais  (Requisito fraco)",
    2
)

if (-not $found) {
    throw "Could not find the LOM3036/LOM3082/LOM3057 requirement block to replace."
}
